# Rename the two enrollment worksheets.
$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General_Enrollment")
$wsGeneral.Name = "Enrollment_General"

$wsRace = $wb.Worksheets.Item("Race_Enrollment")
$wsRace.Name = "Enrollment_Race"

# "High School Units" sheet: widen several columns and move the active
# selection to G11 (it was A19).
$wsHS = $wb.Worksheets.Item("High School Units")

# ColumnWidth is quantized to 1/6-character increments by this engine, so
# these inputs are chosen to land as close as possible to the target
# widths (16.6640625, 13.33203125, 12.6640625, 18.21875, 14.6640625).
$wsHS.Columns.Item(2).ColumnWidth = 15.833333333333334
$wsHS.Columns.Item(3).ColumnWidth = 12.5
$wsHS.Columns.Item(4).ColumnWidth = 11.833333333333334
$wsHS.Columns.Item(6).ColumnWidth = 17.333333333333332
$wsHS.Columns.Item(7).ColumnWidth = 13.833333333333334

# Make "High School Units" the active sheet/tab and move the selection.
$wsHS.Activate() | Out-Null
$wsHS.Range("G11").Select() | Out-Null
